# Scheduled-runner refresh of leve-crafting profit figures (currentAveragePrice*,
# LevePrice*, LeveProfit* columns H:N) across several job sheets.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5050
$ws.Range("I40").Value = 8166.6665
$ws.Range("J40").Value = 3180
$ws.Range("K40").Value = 8166.6665
$ws.Range("L40").Value = 3180
$ws.Range("M40").Value = -7991.6665
$ws.Range("N40").Value = -3530

$ws.Range("H64").Value = 3152.9524
$ws.Range("I64").Value = 3168.3333
$ws.Range("J64").Value = 3146.8
$ws.Range("K64").Value = 3168.3333
$ws.Range("L64").Value = 3146.8
$ws.Range("M64").Value = -2920.3333
$ws.Range("N64").Value = -3642.8

$ws.Range("H67").Value = 3152.9524
$ws.Range("I67").Value = 3168.3333
$ws.Range("J67").Value = 3146.8
$ws.Range("K67").Value = 3168.3333
$ws.Range("L67").Value = 3146.8
$ws.Range("M67").Value = -2310.3333
$ws.Range("N67").Value = -4862.8

$ws.Range("H74").Value = 3348
$ws.Range("I74").Value = 3559.5
$ws.Range("J74").Value = 3277.5
$ws.Range("K74").Value = 3559.5
$ws.Range("L74").Value = 3277.5
$ws.Range("M74").Value = -2623.5
$ws.Range("N74").Value = -5149.5

$ws.Range("H76").Value = 120032.69
$ws.Range("I76").Value = 173647.84
$ws.Range("J76").Value = 3866.5
$ws.Range("K76").Value = 173647.84
$ws.Range("L76").Value = 3866.5
$ws.Range("M76").Value = -173332.84
$ws.Range("N76").Value = -4496.5

$ws.Range("H77").Value = 3348
$ws.Range("I77").Value = 3559.5
$ws.Range("J77").Value = 3277.5
$ws.Range("K77").Value = 17797.5
$ws.Range("L77").Value = 16387.5
$ws.Range("M77").Value = -13117.5
$ws.Range("N77").Value = -25747.5

$ws.Range("H79").Value = 120032.69
$ws.Range("I79").Value = 173647.84
$ws.Range("J79").Value = 3866.5
$ws.Range("K79").Value = 173647.84
$ws.Range("L79").Value = 3866.5
$ws.Range("M79").Value = -172555.84
$ws.Range("N79").Value = -6050.5

$ws.Range("H137").Value = 1367.1522
$ws.Range("I137").Value = 1195.8667
$ws.Range("J137").Value = 1688.3125
$ws.Range("K137").Value = 3587.6001
$ws.Range("L137").Value = 5064.9375
$ws.Range("M137").Value = -1037.6001
$ws.Range("N137").Value = -10164.9375

$ws.Range("H138").Value = 3417
$ws.Range("I138").Value = 3913.5715
$ws.Range("J138").Value = 3352.6296
$ws.Range("K138").Value = 11740.7145
$ws.Range("L138").Value = 10057.8888
$ws.Range("M138").Value = -6600.7145
$ws.Range("N138").Value = -20337.8888

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8834.1
$ws.Range("I32").Value = 4363.875
$ws.Range("J32").Value = 26715
$ws.Range("K32").Value = 4363.875
$ws.Range("L32").Value = 26715
$ws.Range("M32").Value = -4076.875
$ws.Range("N32").Value = -27289

$ws.Range("H63").Value = 2568.9
$ws.Range("I63").Value = 1961.125
$ws.Range("J63").Value = 5000
$ws.Range("K63").Value = 1961.125
$ws.Range("L63").Value = 5000
$ws.Range("M63").Value = -1275.125
$ws.Range("N63").Value = -6372

$ws.Range("H66").Value = 2568.9
$ws.Range("I66").Value = 1961.125
$ws.Range("J66").Value = 5000
$ws.Range("K66").Value = 9805.625
$ws.Range("L66").Value = 25000
$ws.Range("M66").Value = -6373.625
$ws.Range("N66").Value = -31864

$ws.Range("H74").Value = 2554.05
$ws.Range("I74").Value = 2558.9614
$ws.Range("J74").Value = 2544.9285
$ws.Range("K74").Value = 2558.9614
$ws.Range("L74").Value = 2544.9285
$ws.Range("M74").Value = -1684.9614
$ws.Range("N74").Value = -4292.9285

$ws.Range("H77").Value = 2554.05
$ws.Range("I77").Value = 2558.9614
$ws.Range("J77").Value = 2544.9285
$ws.Range("K77").Value = 12794.807
$ws.Range("L77").Value = 12724.6425
$ws.Range("M77").Value = -8426.807000000001
$ws.Range("N77").Value = -21460.6425

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3644.1177
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
# LeveProfitNQ becomes blank (cost == price) rather than a computed loss
$ws.Range("M31").ClearContents()

$ws.Range("H34").Value = 3644.1177
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()

$ws.Range("H62").Value = 45457400
$ws.Range("I62").Value = 2629.2727
$ws.Range("J62").Value = 90912180
$ws.Range("K62").Value = 2629.2727
$ws.Range("L62").Value = 90912180
$ws.Range("M62").Value = -2005.2727
$ws.Range("N62").Value = -90913428

$ws.Range("H65").Value = 45457400
$ws.Range("I65").Value = 2629.2727
$ws.Range("J65").Value = 90912180
$ws.Range("K65").Value = 13146.3635
$ws.Range("L65").Value = 454560900
$ws.Range("M65").Value = -10026.3635
$ws.Range("N65").Value = -454567140

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 839.7632
$ws.Range("I5").Value = 639.5789
$ws.Range("J5").Value = 1039.9474
$ws.Range("K5").Value = 1918.7367
$ws.Range("L5").Value = 3119.8422
$ws.Range("M5").Value = -1806.7367
$ws.Range("N5").Value = -3343.8422

$ws.Range("H131").Value = 985.5
$ws.Range("J131").Value = 1127.238
$ws.Range("L131").Value = 3381.714
$ws.Range("N131").Value = -13461.714

$ws.Range("H135").Value = 839.7632
$ws.Range("I135").Value = 639.5789
$ws.Range("J135").Value = 1039.9474
$ws.Range("K135").Value = 5756.2101
$ws.Range("L135").Value = 9359.526600000001
$ws.Range("M135").Value = -3221.2101
$ws.Range("N135").Value = -14429.5266

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3505.2856
$ws.Range("I80").Value = 6221
$ws.Range("J80").Value = 2656.625
$ws.Range("K80").Value = 6221
$ws.Range("L80").Value = 2656.625
$ws.Range("M80").Value = -5223
$ws.Range("N80").Value = -4652.625

$ws.Range("H83").Value = 3505.2856
$ws.Range("I83").Value = 6221
$ws.Range("J83").Value = 2656.625
$ws.Range("K83").Value = 31105
$ws.Range("L83").Value = 13283.125
$ws.Range("M83").Value = -26113
$ws.Range("N83").Value = -23267.125

$ws.Range("H126").Value = 2589.7856
$ws.Range("I126").Value = 1705.5555
$ws.Range("J126").Value = 4181.4
$ws.Range("K126").Value = 5116.666499999999
$ws.Range("L126").Value = 12544.2
$ws.Range("M126").Value = -2646.666499999999
$ws.Range("N126").Value = -17484.2

$ws.Range("H132").Value = 1966.1
$ws.Range("I132").Value = 1245.8422
$ws.Range("J132").Value = 3210.182
$ws.Range("K132").Value = 3737.5266
$ws.Range("L132").Value = 9630.545999999998
$ws.Range("M132").Value = -1207.5266
$ws.Range("N132").Value = -14690.546
